# Update cryptos list data (prices + 1h volume %) per the Dec 16 2023 refresh,
# and fix the ordering of ARBITRUM / TrustWalletToken / NEARProtocol (rows 47-49).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.427.65'
$ws.Range("E2").Value = '  +0.73%  '
$ws.Range("D3").Value = '2.248.99'
$ws.Range("E3").Value = '  +0.52%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'" + '246.24'
$ws.Range("E5").Value = '  -0.10%  '
$ws.Range("E6").Value = '  +0.32%  '
$ws.Range("D7").Value = "'" + '76.06'
$ws.Range("E7").Value = '  -0.35%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  -1.28%  '
$ws.Range("D10").Value = "'" + '44.91'
$ws.Range("E10").Value = '  +9.29%  '
$ws.Range("D11").Value = "'" + '0.0949'
$ws.Range("E11").Value = '  +0.29%  '
$ws.Range("E12").Value = '  +1.43%  '
$ws.Range("E13").Value = '  -1.27%  '
$ws.Range("D14").Value = "'" + '14.60'
$ws.Range("E14").Value = '  -1.31%  '
$ws.Range("E15").Value = '  -0.02%  '
$ws.Range("D16").Value = '2.254.55'
$ws.Range("E16").Value = '  +0.84%  '
$ws.Range("D17").Value = '42.255.04'
$ws.Range("E17").Value = '  +0.63%  '
$ws.Range("D18").Value = "'" + '0.0000101'
$ws.Range("E18").Value = '  +3.72%  '
$ws.Range("E19").Value = '  +1.05%  '
$ws.Range("D20").Value = "'" + '72.19'
$ws.Range("E20").Value = '  +1.16%  '
$ws.Range("E21").Value = '  +3.46%  '
$ws.Range("D22").Value = "'" + '231.77'
$ws.Range("E22").Value = '  +0.63%  '
$ws.Range("D23").Value = "'" + '8.95'
$ws.Range("E23").Value = '  +23.25%  '
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("D25").Value = "'" + '11.50'
$ws.Range("E25").Value = '  +3.19%  '
$ws.Range("E26").Value = '  -2.18%  '
$ws.Range("E27").Value = '  -0.42%  '
$ws.Range("E28").Value = '  +1.87%  '
$ws.Range("D29").Value = "'" + '167.56'
$ws.Range("E29").Value = '  -0.96%  '
$ws.Range("D30").Value = "'" + '20.66'
$ws.Range("E30").Value = '  +0.80%  '
$ws.Range("D31").Value = "'" + '0.0822'
$ws.Range("E31").Value = '  -3.99%  '
$ws.Range("E32").Value = '  +0.22%  '
$ws.Range("D33").Value = "'" + '30.92'
$ws.Range("E33").Value = '  -6.46%  '
$ws.Range("D34").Value = "'" + '5.29'
$ws.Range("E34").Value = '  +9.85%  '
$ws.Range("E35").Value = '  -0.26%  '
$ws.Range("D36").Value = "'" + '4.56'
$ws.Range("E36").Value = '  -0.79%  '
$ws.Range("D37").Value = "'" + '0.0316'
$ws.Range("E37").Value = '  +6.52%  '
$ws.Range("D38").Value = "'" + '14.04'
$ws.Range("E38").Value = '  +6.45%  '
$ws.Range("E39").Value = '  -0.76%  '
$ws.Range("E40").Value = '  -0.92%  '
$ws.Range("D41").Value = "'" + '63.68'
$ws.Range("E41").Value = '  +6.16%  '
$ws.Range("D42").Value = "'" + '0.201'
$ws.Range("E42").Value = '  -0.53%  '
$ws.Range("D43").Value = "'" + '107.61'
$ws.Range("E43").Value = '  -6.00%  '
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("E45").Value = '  +2.84%  '
$ws.Range("D46").Value = "'" + '0.997'
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = "'" + '2.38'
$ws.Range("E47").Value = '  +5.69%  '
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").Value = "'" + '1.13'
$ws.Range("E48").Value = '  +0.51%  '
$ws.Range("B49").Value = 'TrustWalletToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D49").Value = "'" + '1.19'
$ws.Range("E49").Value = '  +2.48%  '
$ws.Range("E50").Value = '  +0.67%  '
$ws.Range("D51").Value = "'" + '2.70'
$ws.Range("E51").Value = '  +0.79%  '
